# Shop.xlsx edit: insert a new "Force" flag row into the settings block
# (between the existing "Upload" row and the "Desc" header row), pushing
# the item-table rows down by one, and update the frozen-pane / selection
# / data-validation ranges to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row above row 8 ("Upload"), which becomes the new row 8.
$ws.Rows.Item(8).Insert()

# 2. Duplicate the formatting of the row above (row 7, "Ref") into the new
#    row 8 so it keeps the same borders/fills/number formats.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B7:I7").Copy()
$ws.Range("B8:I8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 3. Populate the new row: label + all FALSE boolean flags.
$ws.Range("A8").Value = "Force"
$ws.Range("B8:I8").Value = $false

# 4. Data validations: the ranges below row 7 all grew by one row
#    (the new row 8 needs the same list/blank validation as its
#    neighbours), so rebuild them against the new row numbers.
$ws.Cells.Validation.Delete()
$ws.Range("B7:J9").Validation.Add(3, 1, 1, """TRUE,FALSE""")
$ws.Range("A7:A9").Validation.Add(-4104, 1, 1)

# 5. Frozen pane now covers the extra row (was frozen after row 9, now
#    after row 10), and the remembered selection moves to A9.
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A11").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("A9").Select()
